{"js": "// Fix the typo \"kopelling\" -> \"koppeling\" in the last bullet point\n// (\"Heeft de app een flexibel interface tussen de kopelling van de arduino\n// en de app? ...\").\nconst body = context.document.body;\nconst results = body.search(\"kopelling\", { matchCase: false, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"koppeling\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the typo \"kopelling\" -> \"koppeling\" in the last bullet point\n# (\"Heeft de app een flexibel interface tussen de kopelling van de arduino\n# en de app? ...\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"kopelling\"\n$find.Replacement.Text = \"koppeling\"\n$find.Forward = $true\n$find.Wrap = 1    # wdFindContinue\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
